# Apply the edits described by the diff:
#  - grow the saved workbook window view size (best effort; cosmetic)
#  - append 5 new data rows (315-319) to Sheet1, matching formatting of row 314
#  - set the active selection to D10
#  - dimension / sharedStrings counts update automatically from the new cells

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- try to resize the saved workbook window (cosmetic; may be a no-op in this host) ---
$win = $wb.Windows.Item(1)
$win.Width = 15528
$win.Height = 6648

# --- data for the 5 new rows (columns A..M) ---
# A Color, B Comment, C Document group, D Document name, E Code, F Begin, G End,
# H Weight score, I Segment, J Area, K Coverage %, L Author, M Creation date
$rowNums = 315,316,317,318,319
$colD = "'18945","'10622","'10622","'10622","'137"
$colE = 'Location:Country','Location:Country','Location:City','Location:Hospital name','Location:Country'
$colF = '1: 1660','1: 222','1: 210','1: 156','1: 454'
$colG = '1: 1663','1: 226','1: 214','1: 188','1: 459'
$colI = 'Iran','Italy','Monza','San Gerardo de’ Tintori” Hospital','France'
$colJ = 4,5,5,33,6
$colK = 0.026444532592886424,0.044385264092321346,0.044385264092321346,0.29294274300932088,0.027466239414053559
$colM = '9/17/2019 11:38:55','9/17/2019 11:39:49','9/17/2019 11:40:35','9/17/2019 11:40:57','9/17/2019 11:41:56'

for ($i = 0; $i -lt 5; $i++) {
    $r = $rowNums[$i]

    # 1) set the real values first, left to right, so shared strings are minted in order.
    #    B/C/D are forced to stay text (leading "'" ) because their content looks numeric.
    $ws.Cells.Item($r, 1).Value = '●'
    $ws.Cells.Item($r, 2).Value = "'"
    $ws.Cells.Item($r, 3).Value = "'"
    $ws.Cells.Item($r, 4).Value = $colD[$i]
    $ws.Cells.Item($r, 5).Value = $colE[$i]
    $ws.Cells.Item($r, 6).Value = $colF[$i]
    $ws.Cells.Item($r, 7).Value = $colG[$i]
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = $colI[$i]
    $ws.Cells.Item($r, 10).Value = $colJ[$i]
    $ws.Cells.Item($r, 11).Value = $colK[$i]
    $ws.Cells.Item($r, 12).Value = 'dattaray'
    $ws.Cells.Item($r, 13).Value = $colM[$i]

    # 2) now stamp the row's look (styles/number formats/borders/fills) from the
    #    template row 314 without touching the values we just set.
    $ws.Range("A314:M314").Copy()
    $ws.Range("A" + $r + ":M" + $r).PasteSpecial(-4122)

    # 3) match the template row height (15.6)
    $ws.Rows.Item($r).RowHeight = 15.6
}

# --- set the active selection shown in the saved sheetView ---
[void]$ws.Range("D10").Select()

Write-Output "edit complete"
